$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This sheet lists one "attribute" (player stat) per row: column A holds the
# raw repr captured from the source object, column B extracts the attribute
# name, column C extracts its type. Since a dedicated optimisation class now
# owns identity / team bookkeeping, the raw per-player identifier and
# team-membership attributes are no longer needed here (team constraints will
# be modelled separately), so their rows are dropped.
$namesToRemove = @(
  "'id'",
  "'team_code'",
  "'code'",
  "'squad_number'",
  "'chance_of_playing_this_round'",
  "'chance_of_playing_next_round'",
  "'cost_change_start'",
  "'cost_change_event'",
  "'cost_change_start_fall'",
  "'cost_change_event_fall'",
  "'element_type'",
  "'team'"
)

$lastRow = $ws.UsedRange.Rows.Count
for ($r = $lastRow; $r -ge 2; $r--) {
  $cellText = $ws.Cells.Item($r, 1).Value()
  if ($null -eq $cellText) { continue }
  foreach ($name in $namesToRemove) {
    if ($cellText.StartsWith($name + " (")) {
      $ws.Rows.Item($r).Delete()
      break
    }
  }
}

# Column B's extractor used to leave a trailing space before the "(" marker
# (LEFT stops right before it); wrap it in TRIM so the attribute name is
# clean.
$lastRow = $ws.UsedRange.Rows.Count
$ws.Range("B2:B" + $lastRow).FormulaR1C1 = '=TRIM(LEFT(RC[-1],SEARCH("(",RC[-1])-1))'

$ws.Range("C6").Select() | Out-Null
